# Pedido_Semana_07 - "seguimos con la mejora de los errores de los mails"
#
# Column K ("Stock Mínimo Objetivo") already holds the target minimum-stock
# value per article (rows 3-207). Column L ("Diferencia Stock") was left at
# 0 for every article due to the bug being fixed here; it should mirror the
# Stock Mínimo Objetivo value from column K for each row. The summary cell
# C221 ("Total_Ajuste_Stock") must then reflect the new column-L total,
# matching the already-correct C220 ("Stock_Minimo_Objetivo") total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 3
$lastRow  = 207
$colK     = 11   # K - Stock Mínimo Objetivo
$colL     = 12   # L - Diferencia Stock

$total = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $kVal = $ws.Cells.Item($r, $colK).Value2
    $ws.Cells.Item($r, $colL).Value = $kVal
    $total = $total + $kVal
}

# Total_Ajuste_Stock summary row now matches the sum of the refreshed column L
$ws.Range("C221").Value = $total
